$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 11 match rows that were already in the sheet (rows 2-12) are being
# appended again (rows 13-23) in a different order, exactly reproducing
# the data committed in the source JSON update.
$newRows = @(
    @(" Sharjah", " October 26 2020", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab", "Rahul Tripathi ", "7", "4", "0", "1", "175.00"),
    @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Kolkata Knight Riders", "Sunrisers Hyderabad", "Rahul Tripathi ", "23", "16", "2", "1", "143.75"),
    @(" Dubai (DSC)", " October 29 2020", "Super Kings won by 6 wickets", "Kolkata Knight Riders", "Chennai Super Kings", "Rahul Tripathi ", "3", "2", "0", "0", "150.00"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Kolkata Knight Riders", "Delhi Capitals", "Rahul Tripathi ", "13", "12", "1", "0", "108.33"),
    @(" Abu Dhabi", " October 16 2020", "Mumbai won by 8 wickets (with 19 balls remaining)", "Kolkata Knight Riders", "Mumbai Indians", "Rahul Tripathi ", "7", "9", "1", "0", "77.77"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Rahul Tripathi ", "39", "34", "4", "2", "114.70"),
    @(" Abu Dhabi", " October 10 2020", "KKR won by 2 runs", "Kolkata Knight Riders", "Kings XI Punjab", "Rahul Tripathi ", "4", "10", "1", "0", "40.00"),
    @(" Sharjah", " October 03 2020", "Capitals won by 18 runs", "Kolkata Knight Riders", "Delhi Capitals", "Rahul Tripathi ", "36", "16", "3", "3", "225.00"),
    @(" Abu Dhabi", " October 21 2020", "RCB won by 8 wickets (with 39 balls remaining)", "Kolkata Knight Riders", "Royal Challengers Bangalore", "Rahul Tripathi ", "1", "5", "0", "0", "20.00"),
    @(" Sharjah", " October 12 2020", "RCB won by 82 runs", "Kolkata Knight Riders", "Royal Challengers Bangalore", "Rahul Tripathi ", "16", "22", "1", "0", "72.72"),
    @(" Abu Dhabi", " October 07 2020", "KKR won by 10 runs", "Kolkata Knight Riders", "Chennai Super Kings", "Rahul Tripathi ", "81", "51", "8", "3", "158.82"),
)

$startRow = 13
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Force the whole row to Text format first so that numeric-looking
    # strings (e.g. "175.00", "7") are stored as text, matching the
    # source data (t="str") instead of being coerced into numbers.
    $rowRange = $ws.Range("A$r" + ":K$r")
    $rowRange.NumberFormat = "@"

    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $rowData[$j]
    }
}
